# Weekly update: a new price record is published for "Zapallo italiano" at
# "Macroferia Regional de Talca". The new record is inserted immediately
# above the previous top-of-block data row (row 350), pushing every
# subsequent record down by one row (350-374 -> 351-375).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 350; this shifts rows 350:374 down to
# 351:375 and carries formatting (incl. the date number-format on column D)
# down into the vacated row.
$ws.Rows.Item(350).Insert()

# Populate the newly inserted row 350 with this week's record.
$ws.Cells.Item(350, 1).Value  = 5
$ws.Cells.Item(350, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(350, 3).Value  = "Maule"
$ws.Cells.Item(350, 4).Value  = 44746
$ws.Cells.Item(350, 5).Value  = 7
$ws.Cells.Item(350, 6).Value  = 100112032
$ws.Cells.Item(350, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(350, 8).Value  = "Sin especificar"
$ws.Cells.Item(350, 9).Value  = "Primera"
$ws.Cells.Item(350, 10).Value = 300
$ws.Cells.Item(350, 11).Value = 11000
$ws.Cells.Item(350, 12).Value = 11000
$ws.Cells.Item(350, 13).Value = 11000
$ws.Cells.Item(350, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(350, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(350, 16).Value = 220
$ws.Cells.Item(350, 17).Value = 50
$ws.Cells.Item(350, 18).Value = "Hortaliza"
